$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates that are unambiguously text (coin names, URLs, percent strings,
# and price strings containing two or more dots, e.g. "21.018.67") -- safe to
# assign directly since Excel cannot parse them as a number.
$directTextUpdates = @{
    'D2' = '21.018.67'
    'E2' = '  -4.61%  '
    'D3' = '1.499.71'
    'E3' = '  -3.53%  '
    'E4' = '  +0.55%  '
    'E5' = '  +0.57%  '
    'E6' = '  -2.49%  '
    'E7' = '  -3.76%  '
    'E8' = '  -3.40%  '
    'E9' = '  -2.43%  '
    'E10' = '  -4.97%  '
    'E11' = '  -4.52%  '
    'E12' = '  +0.56%  '
    'E13' = '  -2.62%  '
    'E14' = '  -6.13%  '
    'D15' = '1.506.87'
    'E15' = '  -3.10%  '
    'E16' = '  -4.40%  '
    'E17' = '  -5.87%  '
    'E18' = '  -0.70%  '
    'E19' = '  -2.44%  '
    'E21' = '  -5.51%  '
    'E22' = '  -3.46%  '
    'E23' = '  -4.38%  '
    'E24' = '  -1.76%  '
    'D25' = '21.025.37'
    'E25' = '  -4.65%  '
    'E26' = '  -4.45%  '
    'E27' = '  -1.40%  '
    'E28' = '  -3.87%  '
    'E29' = '  -1.68%  '
    'D30' = '1.674.64'
    'E30' = '  -3.21%  '
    'E31' = '  -4.43%  '
    'E32' = '  -0.97%  '
    'E33' = '  -2.67%  '
    'E34' = '  -5.22%  '
    'E35' = '  -8.25%  '
    'B36' = 'WEMIXTOKEN'
    'C36' = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
    'E36' = '  -7.98%  '
    'B37' = 'InternetComputer(DFINITY)'
    'C37' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'E37' = '  -1.84%  '
    'E38' = '  +2.61%  '
    'E39' = '  -4.24%  '
    'E40' = '  +0.51%  '
    'E41' = '  -6.79%  '
    'E42' = '  -4.66%  '
    'E43' = '  -3.88%  '
    'E44' = '  -4.32%  '
    'E45' = '  -1.58%  '
    'E46' = '  -2.53%  '
    'E47' = '  -3.30%  '
    'B48' = 'EOS'
    'C48' = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
    'E48' = '  -0.80%  '
    'B49' = 'NEARProtocol'
    'C49' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'E49' = '  -3.86%  '
    'E50' = '  -4.85%  '
    'E51' = '  -4.56%  '
}
foreach ($addr in $directTextUpdates.Keys) {
    $ws.Range($addr).Value = $directTextUpdates[$addr]
}

# Price strings that look numeric (e.g. "1.006", "0.5400") must be forced to
# Text format first so Excel keeps them as literal strings (preserving leading
# /trailing zeros and the original decimal text) instead of coercing to a number.
$forcedTextUpdates = @{
    'D4' = '1.006'
    'D6' = '282.96'
    'D7' = '0.3807'
    'D8' = '0.3112'
    'D9' = '42.78'
    'D10' = '0.06898'
    'D11' = '1.031'
    'D13' = '5.543'
    'D14' = '17.67'
    'D16' = '6.346'
    'D17' = '0.00001063'
    'D18' = '0.06541'
    'D19' = '81.62'
    'D21' = '5.937'
    'D22' = '15.01'
    'D23' = '10.82'
    'D24' = '2.334'
    'D26' = '2.313'
    'D27' = '146.63'
    'D28' = '17.85'
    'D29' = '4.792'
    'D31' = '113.60'
    'D32' = '5.784'
    'D33' = '0.9453'
    'D34' = '0.07891'
    'D35' = '8.391'
    'D36' = '1.473'
    'D37' = '5.025'
    'D38' = '10.99'
    'D39' = '0.05744'
    'D40' = '1.005'
    'D41' = '0.02113'
    'D42' = '1.156'
    'D43' = '0.1961'
    'D44' = '0.5574'
    'D45' = '12.89'
    'D46' = '3.651'
    'D47' = '0.5400'
    'D48' = '1.125'
    'D49' = '1.828'
    'D50' = '112.74'
    'D51' = '0.06511'
}
foreach ($addr in $forcedTextUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $forcedTextUpdates[$addr]
    $cell.ClearFormats()
}
